# Fruta / hortaliza, semanal
# The rows 2-19 (data rows) get their per-record fields (Fecha, Volumen,
# Precio minimo/maximo/promedio ponderado, Origen, Precio $/Kg) reshuffled
# across rows - i.e. the same set of weekly records now lines up against
# different row positions. Row 9 keeps its own data (maps to itself).
#
# Mapping: new row R gets the D/M/N/O/P/R/S values that used to live on
# row Source(R), BEFORE any writes happen (so we snapshot everything
# first, then write).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (pre-edit row holding the values to move there)
$map = @{
    2  = 14
    3  = 17
    4  = 8
    5  = 2
    6  = 3
    7  = 10
    8  = 7
    9  = 9
    10 = 15
    11 = 16
    12 = 13
    13 = 6
    14 = 19
    15 = 4
    16 = 5
    17 = 18
    18 = 11
    19 = 12
}

# Columns that actually vary per row and need to move: D,M,N,O,P,R,S
$cols = @(4, 13, 14, 15, 16, 18, 19)   # D, M, N, O, P, R, S

# Snapshot the original values for every row/column we touch before
# writing anything back, so later writes don't clobber values we still
# need to read for other rows.
$snapshot = @{}
foreach ($row in $map.Keys) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Cells.Item($row, $col).Value2
    }
    $snapshot[$row] = $rowVals
}

foreach ($row in $map.Keys) {
    $src = $map[$row]
    $srcVals = $snapshot[$src]
    foreach ($col in $cols) {
        $ws.Cells.Item($row, $col).Value = $srcVals[$col]
    }
}
